$wb = $excel.ActiveWorkbook

# --- Node sheet: add new column G "fixity code" ---------------------------
$ws1 = $wb.Worksheets.Item("Node")

$ws1.Range("G1").Value = "fixity code"

# Rows 2-8: pinned supports -> "RRRFFF"
for ($r = 2; $r -le 8; $r++) {
    $ws1.Cells.Item($r, 7).Value = "RRRFFF"
}

# Rows 9-71: regular nodes -> 0
for ($r = 9; $r -le 71; $r++) {
    $ws1.Cells.Item($r, 7).Value = 0
}

# Rows 72-78: roller supports -> "FRRFFF"
for ($r = 72; $r -le 78; $r++) {
    $ws1.Cells.Item($r, 7).Value = "FRRFFF"
}

# --- Member transformation sheet: trim trailing space in header text ------
$ws4 = $wb.Worksheets.Item("Member transformation")
$ws4.Range("A1").Value = "Transform tag"

# --- Restore / update selections on the non-active sheets ------------------
$ws3 = $wb.Worksheets.Item("Member")
$ws3.Range("B23").Select()

$ws4.Range("B7").Select()

# --- Finally activate Node sheet and set its selection (becomes active tab)
$ws1.Activate()
$ws1.Range("E5").Select()
